$d = $word.ActiveDocument

$pairs = @(
    @("73×42=", "33×67="),
    @("62×22=", "36×93="),
    @("25×50=", "97×11="),
    @("23×30=", "46×24="),
    @("43×78=", "35×42="),
    @("91×82=", "15×15="),
    @("95×72=", "97×26="),
    @("68×38=", "71×37="),
    @("89×46=", "63×19="),
    @("47×50=", "32×26="),
    @("89×41=", "83×47="),
    @("56×81=", "42×29="),
    @("43×57=", "87×85="),
    @("44×25=", "82×50="),
    @("39×67=", "45×81="),
    @("57×22=", "31×24="),
    @("73×64=", "60×11="),
    @("24×47=", "56×72="),
    @("98×65=", "35×38="),
    @("22×23=", "12×38="),
    @("21×62=", "48×50="),
    @("80×56=", "40×39="),
    @("39×53=", "86×15="),
    @("44×51=", "13×20="),
    @("51×38=", "63×37=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
